$d = $word.ActiveDocument

# 1. Update the intro paragraph (whole-run replace keeps xml:space="preserve")
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "B2 Token is a deflationary*") {
        $p.Range.Text = "B2 Token (B2T) is a fully deflationary utility token deployed on BNB Smart Chain. Inspired by the iconic B-2 Spirit, it combines simplicity and strength to build a community-driven ecosystem."
        break
    }
}

# 2. Update the fixed-supply paragraph
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "With a fixed supply of 21 million tokens, B2 is designed*") {
        $p.Range.Text = "With a fixed supply of 21 million tokens and no minting, B2T incentivizes long-term holding and reduces circulating supply over time through a built-in 1% burn on every transfer. There is no transaction tax, ensuring frictionless trading and maximum investor confidence."
        break
    }
}

# 3. Update the Tax bullet value
$d.Content.Find.Execute(
    " 4% on every transfer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " 0% (No transaction tax)",
    2
)

# 4. Remove the "Tax Wallet:" bullet paragraph entirely
foreach ($p in @($d.Paragraphs)) {
    if ($p.Range.Text -like "Tax Wallet:*") {
        $p.Range.Delete()
        break
    }
}
